# Update Set Points in USE_TYPE_PROPERTIES for all UD scenarios
# Sheet "INDOOR_COMFORT" columns B (Tcs_set_C) and C (Ths_set_C) are
# updated for most use-type rows (2-21), shifting the cooling set point
# up by 1 degree and the heating set point down by 1 degree.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("INDOOR_COMFORT")

# New Tcs_set_C (col B) / Ths_set_C (col C) values per row.
$updates = @{
    2  = @(27, 20)   # MULTI_RES
    3  = @(27, 20)   # MULTI_RES_2040
    4  = @(27, 20)   # SINGLE_RES
    5  = @(27, 20)   # HOTEL
    6  = @(27, 20)   # OFFICE
    7  = @(27, 19)   # RETAIL
    8  = @(27, 19)   # FOODSTORE
    9  = @(27, 20)   # RESTAURANT
    10 = @(31, 17)   # INDUSTRIAL
    11 = @(27, 20)   # SCHOOL
    13 = @(27, 17)   # GYM
    14 = @(31, 23)   # SWIMMING
    15 = @(27, 17)   # SERVERROOM
    16 = @(29, 17)   # PARKING
    18 = @(27, 20)   # LAB
    19 = @(27, 20)   # MUSEUM
    20 = @(27, 20)   # LIBRARY
    21 = @(27, 20)   # UNIVERSITY
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
}

# Make INDOOR_COMFORT the active sheet and select B2:C21 (matches the
# saved selection/active-tab state in the workbook after this edit).
$ws.Activate()
$ws.Range("B2:C21").Select()
